$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Version: 0002 -> 0003
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "0003"

# ---------------------------------------------------------------------------
# 2) Descripcion text update (row 7) + new row height
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = 'Se muestran en pantalla todos los datos de la campaña seleccionada incluyendo el botón "Volver".'
$ws.Range("A7:C7").RowHeight = 25.5

# ---------------------------------------------------------------------------
# 3) Insert a new row 11 holding the rich-text "Condicion / Punto de
#    extension" paragraph (between "Puntos de Extension" row and
#    "Curso Basico" row).
# ---------------------------------------------------------------------------
$ws.Rows("11:11").Insert()

# give the new row the same formatting as row 8 (label-blank / text / blank)
$ws.Range("A8:C8").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A11:C11").RowHeight = 51

$ws.Range("B11").Value = 'Condicion: El actor quiere ver el detalle de la campaña. Punto de extensión: Paso 2.c del CU01 Administracion de Supervisores:  El actor hace clic en el botón "Ver" del registro de campaña del formulario "Administración de Supervisores".'

$run1 = $ws.Range("B11").Characters(1, 10)            # "Condicion:"
$run1.Font.Bold = $true

$run2 = $ws.Range("B11").Characters(30, 19)            # "Punto de extensión"
$run2.Font.Bold = $true

# ---------------------------------------------------------------------------
# 4) Curso Basico steps.
#    After the insert above the old rows now sit one lower:
#      row 13 = old step "1. El actor hace clic en Ver..."
#      row 14 = old step "2. El sistema busca..."
#    Target layout:
#      row 13 (step 1) = "El sistema busca..." (re-used)
#      row 14 (step 2) = "El actor hace click en el boton Volver" (new)
#      row 15 (step 3) = "El sistema cierra el formulario." (new row)
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "El sistema busca en la base de datos todos los datos coincidentes con el id de la campaña seleccionada y los presenta en pantalla a través del formulario Campaña."
$ws.Range("B14").Value = 'El actor hace click en el botón "Volver"'

# the long "sistema busca" text used to live on row 14 (it carried a taller
# row height) -- now it lives on row 13, so move the custom height along
# with it and let row 14 fall back to the default height.
$ws.Range("A13:C13").RowHeight = 38.25
$ws.Range("A14:C14").RowHeight = 15

$ws.Rows("15:15").Insert()
$ws.Range("A14:C14").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "El sistema cierra el formulario."

# ---------------------------------------------------------------------------
# 5) Selection ends up on B7, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("B7").Select()

Write-Host "done"
